$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.003420333333333333
$ws.Range("H2").Value = 0.010261
$ws.Range("I2").Value = 0.003549653112303053
$ws.Range("J2").Value = 0.003549653112303053
$ws.Range("M2").Value = 0.668273
$ws.Range("N2").Value = 2.004819
$ws.Range("O2").Value = 0.01328414746766746
$ws.Range("P2").Value = 0.01328414746766746
$ws.Range("Q2").Value = 0.002285716417666667
$ws.Range("R2").Value = 0.020571447759
$ws.Range("S2").Value = 0.00004715411540289852
$ws.Range("T2").Value = 0.00004715411540289851
$ws.Range("G3").Value = 0.003420333333333333
$ws.Range("H3").Value = 0.010261
$ws.Range("I3").Value = 0.003549653112303053
$ws.Range("J3").Value = 0.003549653112303053
$ws.Range("O3").Value = 0.3831531055114357
$ws.Range("P3").Value = 0.3831531055114357
$ws.Range("Q3").Value = 0.06592665023322221
$ws.Range("R3").Value = 0.593339852099
$ws.Range("S3").Value = 0.001360060613467248
$ws.Range("T3").Value = 0.001360060613467248
$ws.Range("G4").Value = 0.003420333333333333
$ws.Range("H4").Value = 0.010261
$ws.Range("I4").Value = 0.003549653112303053
$ws.Range("J4").Value = 0.003549653112303053
$ws.Range("M4").Value = 30.36285833333334
$ws.Range("N4").Value = 91.08857500000001
$ws.Range("O4").Value = 0.6035627470208969
$ws.Range("P4").Value = 0.6035627470208967
$ws.Range("Q4").Value = 0.1038510964527778
$ws.Range("R4").Value = 0.934659868075
$ws.Range("S4").Value = 0.002142438383432907
$ws.Range("T4").Value = 0.002142438383432907
$ws.Range("I5").Value = 0.3907064193682856
$ws.Range("J5").Value = 0.3907064193682855
$ws.Range("M5").Value = 0.668273
$ws.Range("N5").Value = 2.004819
$ws.Range("O5").Value = 0.01328414746766746
$ws.Range("P5").Value = 0.01328414746766746
$ws.Range("Q5").Value = 0.2515862956136667
$ws.Range("R5").Value = 2.264276660523
$ws.Range("S5").Value = 0.005190201691452631
$ws.Range("T5").Value = 0.00519020169145263
$ws.Range("I6").Value = 0.3907064193682856
$ws.Range("J6").Value = 0.3907064193682855
$ws.Range("O6").Value = 0.3831531055114357
$ws.Range("P6").Value = 0.3831531055114357
$ws.Range("S6").Value = 0.149700377924212
$ws.Range("T6").Value = 0.149700377924212
$ws.Range("I7").Value = 0.3907064193682856
$ws.Range("J7").Value = 0.3907064193682855
$ws.Range("M7").Value = 30.36285833333334
$ws.Range("N7").Value = 91.08857500000001
$ws.Range("O7").Value = 0.6035627470208969
$ws.Range("P7").Value = 0.6035627470208967
$ws.Range("S7").Value = 0.235815839752621
$ws.Range("T7").Value = 0.2358158397526209
$ws.Range("I8").Value = 0.6057439275194114
$ws.Range("J8").Value = 0.6057439275194113
$ws.Range("M8").Value = 0.668273
$ws.Range("N8").Value = 2.004819
$ws.Range("O8").Value = 0.01328414746766746
$ws.Range("P8").Value = 0.01328414746766746
$ws.Range("Q8").Value = 0.3900546887903333
$ws.Range("R8").Value = 3.510492199113
$ws.Range("S8").Value = 0.008046791660811931
$ws.Range("T8").Value = 0.008046791660811927
$ws.Range("I9").Value = 0.6057439275194114
$ws.Range("J9").Value = 0.6057439275194113
$ws.Range("O9").Value = 0.3831531055114357
$ws.Range("P9").Value = 0.3831531055114357
$ws.Range("S9").Value = 0.2320926669737565
$ws.Range("T9").Value = 0.2320926669737564
$ws.Range("I10").Value = 0.6057439275194114
$ws.Range("J10").Value = 0.6057439275194113
$ws.Range("M10").Value = 30.36285833333334
$ws.Range("N10").Value = 91.08857500000001
$ws.Range("O10").Value = 0.6035627470208969
$ws.Range("P10").Value = 0.6035627470208967
$ws.Range("S10").Value = 0.365604468884843
$ws.Range("T10").Value = 0.3656044688848429
